$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 41: Hr value changes from 3 to 8
$ws.Range("B41").Value = 8

# Row 41: Task text (D41) gets an appended note about testing/troubleshooting
$ws.Range("D41").Value = "Indie Project: Trying Paula's example function for retrieving an entity based on its own characteristics and that of another entity. Tested it in StoryDaoTest. Used it for the profile display servlet/jsp. Added hibernate assocations for survey tables; reviewed 1:1 relationships. changed DB FK directionality for 1:1 relationships with survey.  Testing and troubleshooting."

# Row 42: A42 gets a new date value (serial 43538 -> 2019-03-14)
$ws.Range("A42").Value = "3/14/2019"

# Row 44: D44 text changes from "plus9:55 - 12:05" to "Thurs 8:50 - x"
$ws.Range("D44").Value = "Thurs 8:50 - x"

# Row 45: D45 cell is removed entirely (was "plus 12:20 - ")
$ws.Range("D45").Clear() | Out-Null

# Scroll/select so the view reflects the new working area, then land the
# active selection on D45 to match the saved selection in the file.
$ws.Activate() | Out-Null
$ws.Range("A35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D45").Select() | Out-Null
